# Update the cryptocurrency price/volume table with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.368.70"
$ws.Range("E2").Value = "  -0.20%  "

$ws.Range("D3").Value = "1.873.55"
$ws.Range("E3").Value = "  -1.19%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9997"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.50"
$ws.Range("E5").Value = "  +0.42%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9994"
$ws.Range("E6").Value = "  +0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4795"
$ws.Range("E7").Value = "  -1.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2820"
$ws.Range("E8").Value = "  -3.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06519"
$ws.Range("E9").Value = "  -1.50%  "

$ws.Range("D10").Value = "1.872.38"
$ws.Range("E10").Value = "  -1.32%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07461"
$ws.Range("E11").Value = "  +1.80%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.57"
$ws.Range("E12").Value = "  -2.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.099"
$ws.Range("E13").Value = "  -1.52%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.22"
$ws.Range("E14").Value = "  +0.26%  "

$ws.Range("E15").Value = "  -1.16%  "

$ws.Range("D16").Value = "30.331.97"
$ws.Range("E16").Value = "  -0.30%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.31"
$ws.Range("E17").Value = "  -1.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007605"
$ws.Range("E19").Value = "  -2.39%  "

$ws.Range("D20").Value = "2.113.62"
$ws.Range("E20").Value = "  -0.55%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.295"
$ws.Range("E21").Value = "  -2.61%  "

$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  +0.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "219.76"
$ws.Range("E23").Value = "  +12.78%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.181"
$ws.Range("E24").Value = "  -0.29%  "

$ws.Range("E25").Value = "  -0.27%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.98"
$ws.Range("E26").Value = "  +1.71%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.46"
$ws.Range("E27").Value = "  +1.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.975"
$ws.Range("E28").Value = "  +1.50%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.449"
$ws.Range("E29").Value = "  -0.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09387"
$ws.Range("E30").Value = "  +2.36%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.317"
$ws.Range("E31").Value = "  +0.12%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.031"
$ws.Range("E32").Value = "  -0.58%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05084"
$ws.Range("E33").Value = "  -0.16%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.203"
$ws.Range("E34").Value = "  +3.75%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7531"
$ws.Range("E35").Value = "  +3.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.712"
$ws.Range("E36").Value = "  +0.68%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01830"
$ws.Range("E37").Value = "  +2.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.611"
$ws.Range("E38").Value = "  -1.53%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.070"
$ws.Range("E39").Value = "  -0.95%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9068"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "106.86"
$ws.Range("E41").Value = "  +0.51%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.907"
$ws.Range("E42").Value = "  +0.57%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4267"
$ws.Range("E43").Value = "  -1.39%  "

$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.394"
$ws.Range("E45").Value = "  -2.10%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.19"
$ws.Range("E46").Value = "  -1.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1284"

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.939"
$ws.Range("E48").Value = "  -0.61%  "

$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.469"
$ws.Range("E49").Value = "  -6.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.66"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3891"
$ws.Range("E51").Value = "  +0.38%  "
